$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full A2:T7 data block (row 2 updated + rows 3-7 appended), per Dr Hou's advice:
# breaking the single FAPs->FAPs self-edge into the full ECs/FAPs/sCs x FAPs/sCs grid.
$data = @(
    @("ECs",  "Agt", "Agtr2", "FAPs", 1, 0.3333333333333333, 0.05671333333333334, 0.17014,    0.02987656731692481, 0.02987656731692481, 3, 1,                  34.682839,           104.048517,  0.9919940127517238,  0.9919940127517238,  1.966979409153333,  17.70281468238,     0.02963737589996325,  0.02963737589996325),
    @("ECs",  "Agt", "Agtr2", "sCs",  1, 0.3333333333333333, 0.05671333333333334, 0.17014,    0.02987656731692481, 0.02987656731692481, 2, 0.6666666666666666, 0.2799113333333333, 0.839734,    0.008005987248276263, 0.008005987248276263, 0.01587470475111111, 0.14287234276,      0.0002391914169615674, 0.0002391914169615674),
    @("FAPs", "Agt", "Agtr2", "FAPs", 3, 1,                  1.761113666666667,  5.283341,   0.9277541615420762,  0.9277541615420761,  3, 1,                  34.682839,           104.048517,  0.9919940127517238,  0.9919940127517238,  61.08042176169967,  549.723795855297,  0.9203265735552352,  0.9203265735552352),
    @("FAPs", "Agt", "Agtr2", "sCs",  3, 1,                  1.761113666666667,  5.283341,   0.9277541615420762,  0.9277541615420761,  2, 0.6666666666666666, 0.2799113333333333, 0.839734,    0.008005987248276263, 0.008005987248276263, 0.4929556745882223,  4.436601071294,    0.007427587986841099, 0.007427587986841098),
    @("sCs",  "Agt", "Agtr2", "FAPs", 1, 0.3333333333333333, 0.08042766666666666, 0.241283,   0.04236927114099899, 0.04236927114099899, 3, 1,                  34.682839,           104.048517,  0.9919940127517238,  0.9919940127517238,  2.789459814145667,  25.105138327311,   0.04203006329652539,  0.04203006329652539),
    @("sCs",  "Agt", "Agtr2", "sCs",  1, 0.3333333333333333, 0.08042766666666666, 0.241283,   0.04236927114099899, 0.04236927114099899, 2, 0.6666666666666666, 0.2799113333333333, 0.839734,    0.008005987248276263, 0.008005987248276263, 0.02251261541355555, 0.202613538722,    0.0003392078444735974, 0.0003392078444735974)
)

$startRow = 2
$colCount = $data[0].Count
for ($c = 0; $c -lt $colCount; $c++) {
    for ($i = 0; $i -lt $data.Count; $i++) {
        $row = $startRow + $i
        $ws.Cells.Item($row, $c + 1).Value = $data[$i][$c]
    }
}
